$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("February 2013")

# Row 28: V28=0, W28=0, X28 copy style+value from Q28 (OFF), Y28=0
$ws.Range("Q28").Copy($ws.Range("X28"))
$ws.Range("V28").Value = 0
$ws.Range("W28").Value = 0
$ws.Range("Y28").Value = 0

# Row 29: V29=0, W29=0, X29 copy style+value from Q29 (OFF), Y29=0
$ws.Range("Q29").Copy($ws.Range("X29"))
$ws.Range("V29").Value = 0
$ws.Range("W29").Value = 0
$ws.Range("Y29").Value = 0

# Row 30: V30=2, W30=1, X30 copy style+value from Q30 (OFF), Y30=2
$ws.Range("Q30").Copy($ws.Range("X30"))
$ws.Range("V30").Value = 2
$ws.Range("W30").Value = 1
$ws.Range("Y30").Value = 2

# Row 31: V31=0, W31=0.5, X31 copy style+value from Q31 (OFF), Y31=2
$ws.Range("Q31").Copy($ws.Range("X31"))
$ws.Range("V31").Value = 0
$ws.Range("W31").Value = 0.5
$ws.Range("Y31").Value = 2

# Update the view: selection moved to Z34
$ws.Range("Z34").Select()

Write-Host "done"
